{"js": "// Update the date line and all 25 three-digit-by-one-digit multiplication\n// equations in the table to the new day's values. Each old value is unique\n// in the document, so we can safely search for the exact old text and\n// replace it in place (which preserves the run/paragraph formatting of the\n// match) rather than rewriting whole paragraphs/cells.\nconst replacements = [\n  [\"2025-01-21 Tuesday\", \"2025-01-22 Wednesday\"],\n  [\"811\u00d72=1622\", \"732\u00d75=3660\"],\n  [\"358\u00d74=1432\", \"932\u00d76=5592\"],\n  [\"134\u00d73=402\", \"583\u00d77=4081\"],\n  [\"185\u00d73=555\", \"520\u00d77=3640\"],\n  [\"468\u00d74=1872\", \"433\u00d77=3031\"],\n  [\"620\u00d79=5580\", \"401\u00d74=1604\"],\n  [\"417\u00d79=3753\", \"894\u00d74=3576\"],\n  [\"157\u00d72=314\", \"933\u00d75=4665\"],\n  [\"518\u00d73=1554\", \"200\u00d75=1000\"],\n  [\"567\u00d74=2268\", \"196\u00d74=784\"],\n  [\"421\u00d78=3368\", \"427\u00d72=854\"],\n  [\"792\u00d75=3960\", \"225\u00d73=675\"],\n  [\"206\u00d75=1030\", \"929\u00d74=3716\"],\n  [\"977\u00d74=3908\", \"953\u00d72=1906\"],\n  [\"675\u00d77=4725\", \"391\u00d75=1955\"],\n  [\"872\u00d75=4360\", \"583\u00d78=4664\"],\n  [\"620\u00d74=2480\", \"980\u00d76=5880\"],\n  [\"434\u00d76=2604\", \"780\u00d74=3120\"],\n  [\"945\u00d74=3780\", \"807\u00d79=7263\"],\n  [\"113\u00d72=226\", \"401\u00d79=3609\"],\n  [\"987\u00d73=2961\", \"806\u00d75=4030\"],\n  [\"227\u00d76=1362\", \"908\u00d78=7264\"],\n  [\"792\u00d79=7128\", \"158\u00d78=1264\"],\n  [\"483\u00d72=966\", \"279\u00d75=1395\"],\n  [\"756\u00d75=3780\", \"350\u00d76=2100\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Not found: \" + oldText);\n  }\n\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and all 25 three-digit-by-one-digit multiplication\n# equations in the table to the new day's values. Each old value is unique\n# in the document, so Find/Replace on the whole document content is safe\n# and keeps each run's original formatting (font/size/alignment) intact.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2025-01-21 Tuesday\", \"2025-01-22 Wednesday\"),\n  @(\"811\u00d72=1622\", \"732\u00d75=3660\"),\n  @(\"358\u00d74=1432\", \"932\u00d76=5592\"),\n  @(\"134\u00d73=402\", \"583\u00d77=4081\"),\n  @(\"185\u00d73=555\", \"520\u00d77=3640\"),\n  @(\"468\u00d74=1872\", \"433\u00d77=3031\"),\n  @(\"620\u00d79=5580\", \"401\u00d74=1604\"),\n  @(\"417\u00d79=3753\", \"894\u00d74=3576\"),\n  @(\"157\u00d72=314\", \"933\u00d75=4665\"),\n  @(\"518\u00d73=1554\", \"200\u00d75=1000\"),\n  @(\"567\u00d74=2268\", \"196\u00d74=784\"),\n  @(\"421\u00d78=3368\", \"427\u00d72=854\"),\n  @(\"792\u00d75=3960\", \"225\u00d73=675\"),\n  @(\"206\u00d75=1030\", \"929\u00d74=3716\"),\n  @(\"977\u00d74=3908\", \"953\u00d72=1906\"),\n  @(\"675\u00d77=4725\", \"391\u00d75=1955\"),\n  @(\"872\u00d75=4360\", \"583\u00d78=4664\"),\n  @(\"620\u00d74=2480\", \"980\u00d76=5880\"),\n  @(\"434\u00d76=2604\", \"780\u00d74=3120\"),\n  @(\"945\u00d74=3780\", \"807\u00d79=7263\"),\n  @(\"113\u00d72=226\", \"401\u00d79=3609\"),\n  @(\"987\u00d73=2961\", \"806\u00d75=4030\"),\n  @(\"227\u00d76=1362\", \"908\u00d78=7264\"),\n  @(\"792\u00d79=7128\", \"158\u00d78=1264\"),\n  @(\"483\u00d72=966\", \"279\u00d75=1395\"),\n  @(\"756\u00d75=3780\", \"350\u00d76=2100\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $find.Execute([ref]$old, $true, $true, $false, $false, $false, $true, 1, $false, [ref]$new, 2) | Out-Null\n}\n"}
